# Infra_Pipeline.pptx edit:
#  - bump the cached "datetimeFigureOut" footer field from 2023-12-07 to
#    2023-12-08 everywhere it is cached (slide master + every slide layout)
#  - rename the "T2" label callout on slide 3 to "T3"

$p = $ppt.ActivePresentation

# --- 1) Slide master date placeholder -------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    $isDate = $false
    try { $isDate = ($sh.PlaceholderFormat.Type -eq 16) } catch { $isDate = $false }
    if ($isDate) {
        $sh.TextFrame.TextRange.Text = "2023-12-08"
    }
}

# --- 2) Every slide layout's date placeholder ------------------------------
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        $isDate = $false
        try { $isDate = ($sh.PlaceholderFormat.Type -eq 16) } catch { $isDate = $false }
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = "2023-12-08"
        }
    }
}

# --- 3) "T2" -> "T3" label on slide 3 --------------------------------------
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "T2") {
            $sh.TextFrame.TextRange.Text = "T3"
        }
    }
}
